# Update the "Förändrad" (changed) date in column C for every data row
# (rows 2-41) of the active sheet from 2023-11-03 (45233) to 2023-11-13 (45243).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 41; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45233) {
        $cell.Value2 = 45243
    }
}
